$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Add 8 new rows to the table (rows 165-172)
for ($i = 0; $i -lt 8; $i++) {
    $lo.ListRows.Add() | Out-Null
}

# Copy formatting (styles) from the last pre-existing row (164) onto the new rows
$ws.Range("B164:G164").Copy()
$ws.Range("B165:G172").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate "abordagem da aula" notes for rows 166-171 first (controls shared-string order)
$ws.Range("F166").Value = "`n4:14`n6. Autenticação e autorização com tokens JWT`n68. Configuração inicial do Spring Security`no simples fato de adicionar as dependencias do JWT ja bloqueia os endpoints`n"
$ws.Range("F167").Value = "`n5:29`n6. Autenticação e autorização com tokens JWT`n68. Configuração inicial do Spring Security`nanotação @EnableWebSecurity do spring security`n"
$ws.Range("F168").Value = "6:14`n6. Autenticação e autorização com tokens JWT`n68. Configuração inicial do Spring Security`nsobrescrever metodo protected void configure (HttpSecurity http)"
$ws.Range("F169").Value = "7:13`n6. Autenticação e autorização com tokens JWT`n68. Configuração inicial do Spring Security`ndefine no metodo sobrescrito que todos os caminhos descritos no vetor (no caso endpoints) serão permitidos o acesso ... para todo os outros, será necessário autenticação"
$ws.Range("F170").Value = "`n9:02`n6. Autenticação e autorização com tokens JWT`n68. Configuração inicial do Spring Security`nconfiguração para o back end De modo geral pode-se desabilitar proteção de ataques a CSRF em sistemas stateless - nosso sistema é stateless, o que significa que nao armazena seção de usuário"
$ws.Range("F171").Value = "9:51`n6. Autenticação e autorização com tokens JWT`n68. Configuração inicial do Spring Security`no teste sugerido na aula não funcionou - para resolver foi necessário incluir um @Bean na classe JacksonConfig para funcionar com o profile de test e com o H2`n@Bean`npublic JavaMailSender jms (){`nreturn new JavaMailSenderImpl();`n}"

# Row 165 "abordagem" (dependencies note)
$ws.Range("F165").Value = "`n3:15 - Adiciona dependencias necessarias para o JWT (Json Web Tokens) funcionar:`n<dependency>`n<groupId>org.springframework.boot</groupId>`n<artifactId>spring-boot-starter-security</artifactId>`n</dependency>`n<dependency>`n<groupId>io.jsonwebtoken</groupId>`n<artifactId>jjwt</artifactId>`n<version>0.7.0</version>`n</dependency>"

# "nome aula" text shared by rows 165-171
$ws.Range("E165").Value = ". Configuração inicial do Spring Security"

# "Nome da Secao" text shared by all new rows
$ws.Range("C165").Value = "Autenticação e autorização com tokens JWT"

# Row 172 abordagem
$ws.Range("F172").Value = "0:57`n6. Autenticação e autorização com tokens JWT`n69. Adicionando senha a Cliente`nCriação de @Bean de BCryptPasswordEncoder no arquivo de configuração - tem a função de criptografar a senha para armazenar no banco de dados"

# Row 172 nome aula
$ws.Range("E172").Value = "Adicionando senha a Cliente"

# Fill remaining shared text (reused strings) for rows 166-171
$ws.Range("E166").Value = ". Configuração inicial do Spring Security"
$ws.Range("C166").Value = "Autenticação e autorização com tokens JWT"
$ws.Range("E167").Value = ". Configuração inicial do Spring Security"
$ws.Range("C167").Value = "Autenticação e autorização com tokens JWT"
$ws.Range("E168").Value = ". Configuração inicial do Spring Security"
$ws.Range("C168").Value = "Autenticação e autorização com tokens JWT"
$ws.Range("E169").Value = ". Configuração inicial do Spring Security"
$ws.Range("C169").Value = "Autenticação e autorização com tokens JWT"
$ws.Range("E170").Value = ". Configuração inicial do Spring Security"
$ws.Range("C170").Value = "Autenticação e autorização com tokens JWT"
$ws.Range("E171").Value = ". Configuração inicial do Spring Security"
$ws.Range("C171").Value = "Autenticação e autorização com tokens JWT"
$ws.Range("C172").Value = "Autenticação e autorização com tokens JWT"

# Fill numeric columns B (Secao) and D (Aula) for all new rows
$ws.Range("B165").Value = 6
$ws.Range("D165").Value = 68
$ws.Range("B166").Value = 6
$ws.Range("D166").Value = 68
$ws.Range("B167").Value = 6
$ws.Range("D167").Value = 68
$ws.Range("B168").Value = 6
$ws.Range("D168").Value = 68
$ws.Range("B169").Value = 6
$ws.Range("D169").Value = 68
$ws.Range("B170").Value = 6
$ws.Range("D170").Value = 68
$ws.Range("B171").Value = 6
$ws.Range("D171").Value = 68
$ws.Range("B172").Value = 6
$ws.Range("D172").Value = 69

# "aprendido" spacer text on row 165 (reuses existing shared string)
$ws.Range("G165").Value = "`n`n`n`n`n`n"

# Row heights to match target layout
$ws.Rows.Item(165).RowHeight = 345
$ws.Rows.Item(166).RowHeight = 105
$ws.Rows.Item(167).RowHeight = 90
$ws.Rows.Item(168).RowHeight = 60
$ws.Rows.Item(169).RowHeight = 90
$ws.Rows.Item(170).RowHeight = 105
$ws.Rows.Item(171).RowHeight = 210
$ws.Rows.Item(172).RowHeight = 90

# Update selection/view to match target state
$ws.Range("C173").Select()
